$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("05-10-2021", 29164, 5357, 6331, -26803),
    @("06-10-2021", 29164, 5357, 6331, -26500),
    @("07-10-2021", 29164, 5357, 6331, -24898),
    @("08-10-2021", 29164, 5357, 6331, -24796),
    @("12-10-2021", 29164, 5357, 6331, -23945)
)

$startRow = 196
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $rowData[0]
    $cellA.ClearFormats()

    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
